$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Constant values shared across all rows in this block (139-156)
$ws.Range("A139:A156").Value = 11
$ws.Range("B139:B156").Value = "Vega Monumental Concepción"
$ws.Range("C139:C156").Value = "Bíobío"
$ws.Range("E139:E156").Value = 8
$ws.Range("F139:F156").Value = 100112023
$ws.Range("G139:G156").Value = "Brócoli"
$ws.Range("H139:H156").Value = "Sin especificar"
$ws.Range("N139:N156").Value = "`$/unidad"
$ws.Range("Q139:Q156").Value = 1
$ws.Range("R139:R156").Value = "Hortaliza"

# Per-row data that differs: Fecha(D), Calidad(I), Volumen(J), Precio minimo(K), Precio maximo(L), Precio promedio ponderado(M), Origen(O), Precio $/Kg(P)
$rows = @(
    @{ Row=139; D=44476; I="Primera"; J=2000; K=600; L=700; M=650; O="Región Metropolitana"; P=650 },
    @{ Row=140; D=44476; I="Segunda"; J=1000; K=500; L=500; M=500; O="Región Metropolitana"; P=500 },
    @{ Row=141; D=44386; I="Primera"; J=1000; K=700; L=800; M=750; O="Provincia de Chacabuco"; P=750 },
    @{ Row=142; D=44386; I="Segunda"; J=500; K=600; L=600; M=600; O="Provincia de Chacabuco"; P=600 },
    @{ Row=143; D=44306; I="Primera"; J=1500; K=600; L=700; M=647; O="Región Metropolitana"; P=647 },
    @{ Row=144; D=44306; I="Segunda"; J=600; K=500; L=500; M=500; O="Región Metropolitana"; P=500 },
    @{ Row=145; D=44357; I="Primera"; J=2000; K=600; L=700; M=650; O="Región Metropolitana"; P=650 },
    @{ Row=146; D=44357; I="Segunda"; J=1000; K=500; L=500; M=500; O="Región Metropolitana"; P=500 },
    @{ Row=147; D=44314; I="Primera"; J=1000; K=700; L=800; M=750; O="Región Metropolitana"; P=750 },
    @{ Row=148; D=44314; I="Segunda"; J=500; K=600; L=600; M=600; O="Región Metropolitana"; P=600 },
    @{ Row=149; D=44425; I="Primera"; J=2000; K=600; L=700; M=650; O="Región Metropolitana"; P=650 },
    @{ Row=150; D=44425; I="Segunda"; J=1000; K=500; L=500; M=500; O="Región Metropolitana"; P=500 },
    @{ Row=151; D=44187; I="Primera"; J=1000; K=600; L=700; M=650; O="Región del Maule"; P=650 },
    @{ Row=152; D=44187; I="Segunda"; J=500; K=500; L=500; M=500; O="Región del Maule"; P=500 },
    @{ Row=153; D=44250; I="Primera"; J=1000; K=900; L=1000; M=950; O="Región Metropolitana"; P=950 },
    @{ Row=154; D=44250; I="Segunda"; J=500; K=800; L=800; M=800; O="Región Metropolitana"; P=800 },
    @{ Row=155; D=44432; I="Primera"; J=2000; K=600; L=700; M=650; O="Región Metropolitana"; P=650 },
    @{ Row=156; D=44432; I="Segunda"; J=1000; K=500; L=500; M=500; O="Región Metropolitana"; P=500 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 9).Value = $r.I
    $ws.Cells.Item($r.Row, 10).Value = $r.J
    $ws.Cells.Item($r.Row, 11).Value = $r.K
    $ws.Cells.Item($r.Row, 12).Value = $r.L
    $ws.Cells.Item($r.Row, 13).Value = $r.M
    $ws.Cells.Item($r.Row, 15).Value = $r.O
    $ws.Cells.Item($r.Row, 16).Value = $r.P
}

# Rows 155-156 are brand new; make sure the Fecha (date) column keeps the same
# date number format used by the rest of the column (applying it to the whole
# block is harmless/idempotent for the rows that already had it).
$ws.Range("D139:D156").NumberFormat = "YYYY-MM-DD HH:MM:SS"
